# Applies the "Added Credentials from Assts" edit:
#  - Settings sheet: fix/extend API request-format rows (OrderUpdate, ClaimUpdate,
#    OrderDownload, ClaimDownload), wrap text + row heights for the JSON bodies.
#  - Constants sheet: add Completed_MailSubject / Completed_MailBody rows.
#  - Assets sheet: add new credential/asset rows pulled in from the APIs &
#    notification settings (Source, Version, Download_LastNDays, mail
#    receivers/copy-to lists, API timeouts).

$wb = $excel.ActiveWorkbook
$wsSettings  = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")
$wsAssets    = $wb.Worksheets.Item("Assets")

# ---------------------------------------------------------------------------
# Settings sheet (sheet1): rows 5-8, Name/Value pairs for API request formats
# ---------------------------------------------------------------------------

$jsonOrderUpdate = @'
{
 "appName": "<appName>",
 "format": "json",
 "param": [{
         "orderId": "<orderId>", "poNbr": "<poNbr>", "receiptNbr":  "<receiptNbr>","resultCode": "<resultCode>","resultMsg": "<resultMsg>"}],
 "sign": "",
 "source": "RPA",
 "timestamp": "<timestamp>",
 "version": "<version>"
}
'@

$jsonClaimUpdate = @'
{
 "appName": "<appName>",
 "format": "json",
 "param": [
         {
    "orderId": "<orderId>",
          "afsId": "<afsId>",
   "compensationId": "<compensationId>",
          "resultCode": "<resultCode>",
          "resultMsg": "<resultMsg>"
         }
         ],
 "sign": "",
 "source": "RPA",
 "timestamp": "<timestamp>",
 "version": "<version>"
}
'@

$jsonDownload = @'
{
 "appName": "<appName>",
 "format": "json",
 "param": {
  "beginTime": "<beginTime>",
  "endTime": "<endTime>"
 },
 "sign": "",
 "source": "<source>",
 "timestamp": "<timestamp>",
 "version": "<version>"
}
'@

$wsSettings.Range("A5").Value = "API_OrderUpdate_RequestFormat"
$wsSettings.Range("B5").Value = $jsonOrderUpdate

$wsSettings.Range("A6").Value = "API_ClaimUpdate_RequestFormat"
$wsSettings.Range("B6").Value = $jsonClaimUpdate

$wsSettings.Range("A7").Value = "API_OrderDownload_RequestFormat"
$wsSettings.Range("B7").Value = $jsonDownload

$wsSettings.Range("A8").Value = "API_ClaimDownload_RequestFormat"
$wsSettings.Range("B8").Value = $jsonDownload

$wsSettings.Range("B5:B8").WrapText = $true

$wsSettings.Rows.Item(5).RowHeight = 165
$wsSettings.Rows.Item(6).RowHeight = 255
$wsSettings.Rows.Item(7).RowHeight = 180
$wsSettings.Rows.Item(8).RowHeight = 180

$wsSettings.Range("B5").Select()

# ---------------------------------------------------------------------------
# Constants sheet (sheet2): new rows 24-25 (Name only, Value left blank)
# ---------------------------------------------------------------------------

$wsConstants.Range("A24").Value = "Completed_MailSubject"
$wsConstants.Range("A25").Value = "Completed_MailBody"

$wsConstants.Range("A24").Select()

# ---------------------------------------------------------------------------
# Assets sheet (sheet3): new rows 11-21, Name/Value pairs
# ---------------------------------------------------------------------------

$assetRows = @(
    @("API_Source", "1128_API_Source"),
    @("API_Version", "1128_API_Version"),
    @("API_Download_LastNDays", "1128_API_Download_LastNDays"),
    @("BusinessException_MailReceiver", "1128_BusinessException_MailReceiver"),
    @("SystemException_MailReceiver", "1128_SystemException_MailReceiver"),
    @("Notification_MailReceiver", "1128_Notification_MailReceiver"),
    @("BusinessException_MailCopyTo", "1128_BusinessException_MailCopyTo"),
    @("SystemException_MailCopyTo", "1128_SystemException_MailCopyTo"),
    @("Notification_MailCopyTo", "1128_Notification_MailCopyTo"),
    @("API_UpdateTimeOut_InSeconds", "1128_API_UpdateTimeOut_InSeconds"),
    @("API_DownloadTimeOut_InSeconds", "1128_API_DownloadTimeOut_InSeconds")
)

$row = 11
foreach ($pair in $assetRows) {
    $wsAssets.Cells.Item($row, 1).Value = $pair[0]
    $wsAssets.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

$wsAssets.Range("A18").Select()
